$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5574305
$ws.Range("J40").Value = 8335833
$ws.Range("L40").Value = 8335833
$ws.Range("N40").Value = -8336183
$ws.Range("H43").Value = 4700.5
$ws.Range("I43").Value = 4700.5
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 4700.5
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -4631.5
$ws.Range("N43").ClearContents()
$ws.Range("H139").Value = 109999
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 109999
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 109999
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -120279
$ws.Range("H140").Value = 59899.5
$ws.Range("J140").Value = 59899.5
$ws.Range("L140").Value = 59899.5
$ws.Range("N140").Value = -70259.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 12
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 3293855.8
$ws.Range("I32").Value = 3680663.2
$ws.Range("J32").Value = 5991.5
$ws.Range("K32").Value = 3680663.2
$ws.Range("L32").Value = 5991.5
$ws.Range("M32").Value = -3680376.2
$ws.Range("N32").Value = -6565.5
$ws.Range("H45").Value = 4240.6113
$ws.Range("I45").Value = 1605.3572
$ws.Range("J45").Value = 13464
$ws.Range("K45").Value = 1605.3572
$ws.Range("L45").Value = 13464
$ws.Range("M45").Value = -1228.3572
$ws.Range("N45").Value = -14218
$ws.Range("H74").Value = 31670.377
$ws.Range("I74").Value = 36473.066
$ws.Range("K74").Value = 36473.066
$ws.Range("M74").Value = -35599.066
$ws.Range("H77").Value = 31670.377
$ws.Range("I77").Value = 36473.066
$ws.Range("K77").Value = 182365.33
$ws.Range("M77").Value = -177997.33
$ws.Range("H122").Value = 9983.034
$ws.Range("I122").Value = 10380.56
$ws.Range("K122").Value = 31141.68
$ws.Range("M122").Value = -28691.68
$ws.Range("H132").Value = 9823.5
$ws.Range("I132").Value = 9264.333000000001
$ws.Range("J132").Value = 10468.692
$ws.Range("K132").Value = 27792.999
$ws.Range("L132").Value = 31406.076
$ws.Range("M132").Value = -25262.999
$ws.Range("N132").Value = -36466.076
$ws.Range("H133").Value = 87505.125
$ws.Range("J133").Value = 87505.125
$ws.Range("L133").Value = 87505.125
$ws.Range("N133").Value = -92565.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 15
$ws.Range("N4").ClearContents()
$ws.Range("H109").Value = 59382.332
$ws.Range("J109").Value = 59382.332
$ws.Range("L109").Value = 59382.332
$ws.Range("N109").Value = -62156.332
$ws.Range("H134").Value = 4533.754
$ws.Range("I134").Value = 2333.8125
$ws.Range("K134").Value = 7001.4375
$ws.Range("M134").Value = -4466.4375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1287.2
$ws.Range("I7").Value = 5050
$ws.Range("J7").Value = 346.5
$ws.Range("K7").Value = 5050
$ws.Range("L7").Value = 346.5
$ws.Range("M7").Value = -4937
$ws.Range("N7").Value = -572.5
$ws.Range("H31").Value = 6323.9844
$ws.Range("I31").Value = 2669.9
$ws.Range("J31").Value = 12170.52
$ws.Range("K31").Value = 2669.9
$ws.Range("L31").Value = 12170.52
$ws.Range("M31").Value = -2374.9
$ws.Range("N31").Value = -12760.52
$ws.Range("H34").Value = 6323.9844
$ws.Range("I34").Value = 2669.9
$ws.Range("J34").Value = 12170.52
$ws.Range("K34").Value = 2669.9
$ws.Range("L34").Value = 12170.52
$ws.Range("M34").Value = -2467.9
$ws.Range("N34").Value = -12574.52
$ws.Range("H122").Value = 1446.3572
$ws.Range("I122").Value = 1108.375
$ws.Range("J122").Value = 1897
$ws.Range("K122").Value = 3325.125
$ws.Range("L122").Value = 5691
$ws.Range("M122").Value = -875.125
$ws.Range("N122").Value = -10591

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5002874.5
$ws.Range("I5").Value = 13333800
$ws.Range("K5").Value = 40001400
$ws.Range("M5").Value = -40001288
$ws.Range("H68").Value = 3728.25
$ws.Range("I68").Value = 1140
$ws.Range("K68").Value = 3420
$ws.Range("M68").Value = -2609
$ws.Range("H71").Value = 3728.25
$ws.Range("I71").Value = 1140
$ws.Range("K71").Value = 10260
$ws.Range("M71").Value = -6204
$ws.Range("H75").Value = 30304896
$ws.Range("I75").Value = 83334040
$ws.Range("J75").Value = 18520642
$ws.Range("K75").Value = 250002120
$ws.Range("L75").Value = 55561926
$ws.Range("M75").Value = -250001122
$ws.Range("N75").Value = -55563922
$ws.Range("H78").Value = 30304896
$ws.Range("I78").Value = 83334040
$ws.Range("J78").Value = 18520642
$ws.Range("K78").Value = 750006360
$ws.Range("L78").Value = 166685778
$ws.Range("M78").Value = -750001368
$ws.Range("N78").Value = -166695762
$ws.Range("H122").Value = 2831832
$ws.Range("I122").Value = 5658725.5
$ws.Range("J122").Value = 4938.8
$ws.Range("K122").Value = 50928529.5
$ws.Range("L122").Value = 44449.2
$ws.Range("M122").Value = -50926079.5
$ws.Range("N122").Value = -49349.2
$ws.Range("H135").Value = 5002874.5
$ws.Range("I135").Value = 13333800
$ws.Range("K135").Value = 120004200
$ws.Range("M135").Value = -120001665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1974.7693
$ws.Range("I97").Value = 1778.909
$ws.Range("J97").Value = 2228.2354
$ws.Range("K97").Value = 1778.909
$ws.Range("L97").Value = 2228.2354
$ws.Range("M97").Value = -1282.909
$ws.Range("N97").Value = -3220.2354
$ws.Range("H99").Value = 9568.5
$ws.Range("I99").Value = 8482.200000000001
$ws.Range("J99").Value = 15000
$ws.Range("K99").Value = 8482.200000000001
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = -6236.200000000001
$ws.Range("N99").Value = -19492
$ws.Range("H107").Value = 1342.4286
$ws.Range("I107").Value = 1224.75
$ws.Range("J107").Value = 1499.3334
$ws.Range("K107").Value = 1224.75
$ws.Range("L107").Value = 1499.3334
$ws.Range("M107").Value = 695.25
$ws.Range("N107").Value = -5339.3334
$ws.Range("H132").Value = 2081
$ws.Range("I132").Value = 2035.05
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 6105.15
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -3575.15
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3982
$ws.Range("I22").Value = 1999
$ws.Range("J22").Value = 4477.75
$ws.Range("K22").Value = 1999
$ws.Range("L22").Value = 4477.75
$ws.Range("M22").Value = -1704
$ws.Range("N22").Value = -5067.75
$ws.Range("H27").Value = 3982
$ws.Range("I27").Value = 1999
$ws.Range("J27").Value = 4477.75
$ws.Range("K27").Value = 1999
$ws.Range("L27").Value = 4477.75
$ws.Range("M27").Value = -1892
$ws.Range("N27").Value = -4691.75
$ws.Range("H61").Value = 5806.857
$ws.Range("I61").Value = 4749.3335
$ws.Range("J61").Value = 6600
$ws.Range("K61").Value = 4749.3335
$ws.Range("L61").Value = 6600
$ws.Range("M61").Value = -4547.3335
$ws.Range("N61").Value = -7004
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774
$ws.Range("H113").Value = 5806.857
$ws.Range("I113").Value = 4749.3335
$ws.Range("J113").Value = 6600
$ws.Range("K113").Value = 4749.3335
$ws.Range("L113").Value = 6600
$ws.Range("M113").Value = -2579.3335
$ws.Range("N113").Value = -10940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1474.4615
$ws.Range("I107").Value = 1519.4445
$ws.Range("J107").Value = 1373.25
$ws.Range("K107").Value = 4558.333500000001
$ws.Range("L107").Value = 4119.75
$ws.Range("M107").Value = -2638.333500000001
$ws.Range("N107").Value = -7959.75
$ws.Range("H126").Value = 2544.1936
$ws.Range("I126").Value = 994.5833
$ws.Range("K126").Value = 2983.7499
$ws.Range("M126").Value = -513.7498999999998
$ws.Range("H132").Value = 41740036
$ws.Range("I132").Value = 62516056
$ws.Range("J132").Value = 187999.75
$ws.Range("K132").Value = 187548168
$ws.Range("L132").Value = 563999.25
$ws.Range("M132").Value = -187545638
$ws.Range("N132").Value = -569059.25
$ws.Range("H141").Value = 67465
$ws.Range("J141").Value = 67465
$ws.Range("L141").Value = 67465
$ws.Range("N141").Value = -77825
